$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card22")

# Row 18: fill previously-empty B18:K18 with the literal text "nan"
$ws.Range("B18:K18").Value = "nan"

# Row 19: new event row
# Column A holds text "22" (matches the existing text-typed "card" id in column A)
$ws.Range("A19").Value = "'22"
$ws.Range("A19").Style = "Normal"

# Columns B:K are present but empty (text cells with no content), mirroring row 18's prior state
$ws.Range("B19:K19").Value = "'"
$ws.Range("B19:K19").Style = "Normal"

$ws.Range("L19").Value = "14\8\2025"
$ws.Range("M19").Value = "9998 h"
$ws.Range("N19").Value = "تم تغيير زيت الجيربوكس"
$ws.Range("O19").Value = "تيم العمل"
